$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.439.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.66%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.640.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.64%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.32%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.9986"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'305.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.78%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.3614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -5.54%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'47.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -5.30%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.3243"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -9.90%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -9.40%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.06889"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -10.34%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.9978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.26%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'5.905"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -9.15%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'19.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -11.32%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'1.635.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.61%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'6.510"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -8.80%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001041"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -9.67%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.06496"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.23%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.05%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'76.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -11.37%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'5.873"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -10.03%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'15.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -11.60%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'12.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -8.42%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'24.390.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.68%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.391"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.33%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.324"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -20.26%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'145.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -6.25%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'18.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -11.13%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.817.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.62%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'123.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -7.60%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.141"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").Value = "'4.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.81%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.571"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -23.01%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.08318"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.05%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.662"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -7.63%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'12.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -14.17%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'5.129"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -10.87%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.06005"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -10.98%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.02206"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -11.72%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'1.199"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -7.51%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'8.175"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -12.67%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.2023"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -10.60%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.13%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.5824"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -11.58%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'3.713"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'12.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -13.33%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -12.22%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'121.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -7.88%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.923"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -11.99%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.06883"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -7.96%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'73.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -9.54%  "
$ws.Range("E51").Style = "Normal"
